# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# worksheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$ws1.Range("F2").Value  = 124
$ws1.Range("F4").Value  = 933
$ws1.Range("F5").Value  = 974
$ws1.Range("F6").Value  = 1738
$ws1.Range("F7").Value  = 384
$ws1.Range("F8").Value  = 1159
$ws1.Range("F9").Value  = 50
$ws1.Range("F11").Value = 115
$ws1.Range("F12").Value = 267
$ws1.Range("F13").Value = 47
$ws1.Range("F16").Value = 139
$ws1.Range("F17").Value = 91
$ws1.Range("F21").Value = 110
$ws1.Range("F23").Value = 17
$ws1.Range("F24").Value = 632
$ws1.Range("F25").Value = 139
$ws1.Range("F27").Value = 848
$ws1.Range("F29").Value = 124
$ws1.Range("F30").Value = 28
$ws1.Range("F31").Value = 253
$ws1.Range("F33").Value = 12

# --- 全部类型 (sheet4) ---
$ws4.Range("F3").Value  = 124
$ws4.Range("F5").Value  = 933
$ws4.Range("F6").Value  = 974
$ws4.Range("F7").Value  = 1738
$ws4.Range("F8").Value  = 384
$ws4.Range("F9").Value  = 1159
$ws4.Range("F10").Value = 50
$ws4.Range("F13").Value = 115
$ws4.Range("F14").Value = 267
$ws4.Range("F15").Value = 47
$ws4.Range("F18").Value = 139
$ws4.Range("F19").Value = 91
$ws4.Range("F29").Value = 110
$ws4.Range("F31").Value = 17
$ws4.Range("F32").Value = 632
$ws4.Range("F33").Value = 139
$ws4.Range("F35").Value = 848
$ws4.Range("F39").Value = 124
$ws4.Range("F40").Value = 28
$ws4.Range("F41").Value = 253
$ws4.Range("F46").Value = 12
